## Rename the embedded logo pictures' shape names:
##   - Pearson logo inline picture in the "first page" footer  (id=3) : image1.png -> image2.png
##   - Pearson logo inline picture in the "default" footer     (id=2) : image1.png -> image2.png
##   - BTec logo inline picture in the "first page" header     (id=1) : image2.jpg -> image1.jpg
##
## InlineShape has no settable/gettable .Name on the Word object model (same as
## real Word) - the only way to rename the underlying drawing object is to
## temporarily convert the inline picture to a floating shape (which *does*
## expose .Name), set the new name, then convert it back to an inline shape so
## the <wp:inline> wrapper (and everything else) is preserved.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlinePicture($range, $newName) {
    $count = $range.InlineShapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $inline = $range.InlineShapes.Item($i)
        $shape = $inline.ConvertToShape()
        $shape.Name = $newName
        [void]$shape.ConvertToInlineShape()
    }
}

# Footer, first page (docPr id="3") - Pearson logo: image1.png -> image2.png
Rename-InlinePicture $sec.Footers.Item(2).Range "image2.png"

# Footer, default/odd pages (docPr id="2") - Pearson logo: image1.png -> image2.png
Rename-InlinePicture $sec.Footers.Item(1).Range "image2.png"

# Header, first page (docPr id="1") - BTec logo: image2.jpg -> image1.jpg
Rename-InlinePicture $sec.Headers.Item(2).Range "image1.jpg"
